# Apply the "single child" update: the old "max" column (C) is dropped,
# the "prediction" values (old column B) are replaced with real numeric
# prediction scores, and the remaining columns (prediction / rejection-f)
# shift left to fill the gap left by the removed "max" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column B ("1-f__UBA660") with the new numeric prediction values
#    while the sheet still has its original column layout (A..E).
$ws.Cells.Item(2, 2).Value = 9400.745327406312
$ws.Cells.Item(3, 2).Value = 17266.73476994588
$ws.Cells.Item(4, 2).Value = 31017.43220179664
$ws.Cells.Item(5, 2).Value = 8532.149193212354
$ws.Cells.Item(6, 2).Value = 13344.36234376533
$ws.Cells.Item(7, 2).Value = 2553.426947584233
$ws.Cells.Item(8, 2).Value = 4886.493564873062
$ws.Cells.Item(9, 2).Value = 15270.18987267277
$ws.Cells.Item(10, 2).Value = 9075.031332159211
$ws.Cells.Item(11, 2).Value = 10246.47842293522
$ws.Cells.Item(12, 2).Value = 22227.07057721603
$ws.Cells.Item(13, 2).Value = 17133.23427466815
$ws.Cells.Item(14, 2).Value = 26992.90070174724
$ws.Cells.Item(15, 2).Value = 18731.50036110093
$ws.Cells.Item(16, 2).Value = 5057.610322543561
$ws.Cells.Item(17, 2).Value = 24983.16456973581
$ws.Cells.Item(18, 2).Value = 17470.46482815067
$ws.Cells.Item(19, 2).Value = 20590.44830364989
$ws.Cells.Item(20, 2).Value = 21022.35645824457
$ws.Cells.Item(21, 2).Value = 18835.15125770698
$ws.Cells.Item(22, 2).Value = 8444.917494803041
$ws.Cells.Item(23, 2).Value = 15611.13397791301
$ws.Cells.Item(24, 2).Value = 17127.09134513855

# 2) Remove the obsolete "max" column (C), shifting the "prediction" and
#    "rejection-f" columns one place to the left so that D24 becomes the
#    new rightmost used cell (A1:D24).
$ws.Range("C:C").Delete()

# 3) Refresh the worksheet's used range/dimension to reflect the new
#    A1:D24 extent.
$ws.UsedRange | Out-Null
